# Update "想去人数" (want-to-go headcount) figures in column F for the
# "展览" (Exhibition) sheet and the aggregated "全部类型" (All Types) sheet.
# Same events are listed in both sheets (at different row numbers), so each
# sheet gets its own row -> new-value map.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> new F value, for the "展览" sheet.
$exhibitionUpdates = @{
    3  = 11014
    5  = 1225
    6  = 1101
    7  = 851
    8  = 289
    10 = 1179
    12 = 154
    13 = 900
    15 = 2051
    17 = 1018
    18 = 841
    21 = 928
    24 = 642
    25 = 668
    28 = 1026
    31 = 180
    35 = 1903
    36 = 397
    38 = 1453
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new F value, for the "全部类型" sheet.
$allTypesUpdates = @{
    5  = 11014
    9  = 1101
    11 = 1179
    13 = 154
    14 = 900
    15 = 2051
    17 = 1018
    18 = 841
    21 = 928
    25 = 642
    28 = 668
    31 = 1026
    35 = 180
    39 = 1453
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
